# Update the "Use x0y0_direct ..." step label in the workflow diagram to
# also mention x0y0_to_plane, and grow the textbox to fit the extra line
# of wrapped text (matches the author's spAutoFit behaviour in real
# PowerPoint).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item("TextBox 140")

# Update the text first so any auto-fit sizing is based on the new text...
$shape.TextFrame.TextRange.Text = "Use x0y0_direct or x0y0_to_plane to solve magnetic field streamlines"

# ...then pin the height to the exact value from the canonical edit
# (601511 EMU -> points, since Shape.Height is expressed in points).
$shape.Height = 601511 / 12700
